# cypress/TestGrid.xlsx -- "feat: updated grid with newly working test"
#
# The "Branching" feature's VerifyEditTrainingControlsAndLabels turn got a
# newly-working automated test, so several rows in the grid that previously
# had no "Test Name" (column E) now reference the "Branching" test, and the
# row that used to be marked "manual" is now covered by that same test too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32-34 and 36 (Branching feature, Train Dialog page) previously had no
# Test Name in column E -- they are now covered by the "Branching" test.
$ws.Range("E32").Value = "Branching"
$ws.Range("E33").Value = "Branching"
$ws.Range("E34").Value = "Branching"
$ws.Range("E36").Value = "Branching"

# Row 37 was marked "manual" -- it's now automated by the same "Branching" test.
$ws.Range("E37").Value = "Branching"

# Reflect where the author was last working in the grid when the file was saved.
[void]$ws.Range("D35").Select()
